$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("G2").Value = "LEI"
$ws.Range("K2").Value = "['Intocmire lista cantitati - Pasaj superior peste CF Port Constanța']"
$ws.Range("L2").Value = 83.75
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 9

# Delete row 3 entirely (entire row)
$ws.Rows("3").Delete()
